$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix wording of the existing "menu" column header (C1): "du menu" -> "des menu"
$ws.Range("C1").Value = "Création et `ngestion des menu"

# --- Change classe / reassignment of existing task owners ---
# Gestion Invertaire (K) was done by "Matt & Yoan" -> now just "Matt "
$ws.Range("K2").Value = "Matt "

# Gestion Skill (M) header gets extended and the doer changes from Matt to Yoan
$ws.Range("M1").Value = "Gestion Skill et mana"
$ws.Range("M2").Value = "Yoan "

# --- Add new "armes" related columns (O, P) ---
$ws.Range("O1").Value = "Gestion Experience"
$ws.Range("O2").Value = "Enzo"

$ws.Range("P1").Value = "Gestion erreur "

# Match the wrapped header style used by the rest of row 1
$ws.Range("O1:P1").WrapText = $true

# Leave the selection on the newly added last cell, like in the saved workbook
[void]$ws.Range("P2").Select()
